$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill C13:G13 with 5 (style of these cells is already set, so it is preserved)
$ws.Range("C13:G13").Value = 5

# Move the active selection to G13, matching the saved view state in the diff
$ws.Range("G13").Select()
